# Update "想去人数" (F column) figures on the "展览" (sheet1), "演出" (sheet2)
# and "全部类型" (sheet4) worksheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$ws1.Range("F2").Value = 14947
$ws1.Range("F3").Value = 18776
$ws1.Range("F5").Value = 131
$ws1.Range("F13").Value = 55
$ws1.Range("F14").Value = 128
$ws1.Range("F15").Value = 212
$ws1.Range("F16").Value = 59
$ws1.Range("F17").Value = 1437
$ws1.Range("F20").Value = 91
$ws1.Range("F21").Value = 234
$ws1.Range("F22").Value = 7809
$ws1.Range("F24").Value = 29
$ws1.Range("F25").Value = 58
$ws1.Range("F26").Value = 1231
$ws1.Range("F27").Value = 16
$ws1.Range("F28").Value = 5991
$ws1.Range("F32").Value = 154
$ws1.Range("F33").Value = 273
$ws1.Range("F34").Value = 5360

# --- 演出 sheet ---
$ws2.Range("F3").Value = 11

# --- 全部类型 sheet ---
$ws4.Range("F2").Value = 14947
$ws4.Range("F3").Value = 18776
$ws4.Range("F5").Value = 131
$ws4.Range("F13").Value = 55
$ws4.Range("F14").Value = 128
$ws4.Range("F15").Value = 212
$ws4.Range("F16").Value = 59
$ws4.Range("F17").Value = 1437
$ws4.Range("F21").Value = 91
$ws4.Range("F22").Value = 234
$ws4.Range("F23").Value = 7809
$ws4.Range("F25").Value = 29
$ws4.Range("F26").Value = 58
$ws4.Range("F27").Value = 1231
$ws4.Range("F28").Value = 16
$ws4.Range("F29").Value = 11
$ws4.Range("F31").Value = 5991
$ws4.Range("F35").Value = 154
$ws4.Range("F36").Value = 273
$ws4.Range("F37").Value = 5360

$wb.Save()
